$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.464.56'
$ws.Range('E2').Value = '  +3.34%  '
$ws.Range('D3').Value = '3.067.87'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'549.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').Value = "'140.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.062.00'
$ws.Range('E8').Value = '  +2.01%  '
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('D10').Value = "'6.53"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.30%  '
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').Value = "'34.87"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').Value = '3.562.33'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D16').Value = '63.387.83'
$ws.Range('E16').Value = '  +3.24%  '
$ws.Range('D17').Value = '3.066.21'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('E19').Value = '  +1.54%  '
$ws.Range('D20').Value = "'481.93"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').Value = "'13.68"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.99%  '
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('D23').Value = "'7.28"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.20%  '
$ws.Range('D24').Value = "'80.78"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').Value = "'12.64"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.44%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = "'2.76"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').Value = "'7.93"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').Value = "'2.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.90%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = "'26.13"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').Value = "'55.58"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  +1.19%  '
$ws.Range('D37').Value = "'468.02"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('D38').Value = "'0.0821"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('D39').Value = "'0.0397"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.63%  '
$ws.Range('D40').Value = '3.070.90'
$ws.Range('E40').Value = '  -4.21%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('E43').Value = '  +2.75%  '
$ws.Range('D44').Value = "'27.97"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('E45').Value = '  +3.12%  '
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').Value = "'116.48"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.02%  '
$ws.Range('D50').Value = '0.0₃0510'
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('E51').Value = '  +3.02%  '
